$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.58"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.72"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.464"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05761"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.417"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.308"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8194"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.045"
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1420"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07269"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03163"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03121"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.108"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09359"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001607"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04791"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006222"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004132"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009875"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001490"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.727"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1299"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003989"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03846"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006638"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1055"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002670"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006503"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005595"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3890"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002096"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01008"
